$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data row (row 4) - values are assigned in this specific order so that
# the shared-string table grows in the same sequence as the source edit.
$ws.Range("C4").Value = "全球范围"
$ws.Range("D4").Value = "2017-2022"
$ws.Range("E4").Value = "16天"
$ws.Range("F4").Value = "2.25 km x 1.29 km"
$ws.Range("B4").Value = "OCO-2"
$ws.Range("H4").Value = "h5"
$ws.Range("I4").Value = "https://disc.gsfc.nasa.gov/datasets/OCO2_Eph_11r/summary"
$ws.Range("G4").Value = "未知"
$ws.Range("J4").Value = "Processing Level为0，Eph对于数据定位和校正很关键，Att对于数据处理和分析很重要"

# F4 gets a dedicated small grey font (sz 9, #333333, 宋体) - matches the
# new font added to the workbook's style table.
$ws.Range("F4").Font.Size = 9
$ws.Range("F4").Font.Color = 3355443

# I4 becomes a live hyperlink pointing at the NASA dataset page, styled the
# same way as the existing hyperlink cell I3.
$ws.Hyperlinks.Add($ws.Range("I4"), "https://disc.gsfc.nasa.gov/datasets/OCO2_Eph_11r/summary")
$ws.Range("I3").Copy()
$ws.Range("I4").PasteSpecial(-4122)

# Match the author's final selection state.
$ws.Range("C7").Select()
